$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.275.23"
$ws.Range("E2").Value = "  -0.56%  "
$ws.Range("D3").Value = "3.088.04"
$ws.Range("E3").Value = "  +2.40%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "388.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.54"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.538"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.28%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -1.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.98"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.71%  "
$ws.Range("E11").Value = "  +0.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0855"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.69%  "
$ws.Range("D13").Value = "3.580.75"
$ws.Range("E13").Value = "  +2.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.44"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.14%  "
$ws.Range("E15").Value = "  +0.44%  "
$ws.Range("D16").Value = "3.091.57"
$ws.Range("E16").Value = "  +2.42%  "
$ws.Range("E17").Value = "  +2.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.64"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.83%  "
$ws.Range("D19").Value = "51.358.44"
$ws.Range("E19").Value = "  -0.32%  "
$ws.Range("E20").Value = "  +5.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.45"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.26%  "
$ws.Range("D22").Value = "0.0₃0962"
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.02"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "265.54"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.35%  "
$ws.Range("E26").Value = "  -3.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.91%  "
$ws.Range("E28").Value = "  -4.27%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("E30").Value = "  -4.03%  "
$ws.Range("E31").Value = "  -1.94%  "
$ws.Range("E32").Value = "  +2.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "36.02"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.72%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0474"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.20%  "
$ws.Range("E35").Value = "  +0.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "49.78"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.42%  "
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.41"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.290"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "129.35"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.13%  "
$ws.Range("B41").Value = "NEARProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.86"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.81%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.84"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.52%  "
$ws.Range("E43").Value = "  -1.71%  "
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.49"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "21.98"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.81%  "
$ws.Range("E47").Value = "  +5.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.08"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.53%  "
$ws.Range("D49").Value = "2.071.47"
$ws.Range("E49").Value = "  +2.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.936"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +19.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0324"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.72%  "
